$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2427
$ws.Range("I94").Value = 2441.111
$ws.Range("J94").Value = 2300
$ws.Range("K94").Value = 2441.111
$ws.Range("L94").Value = 2300
$ws.Range("M94").Value = -1990.111
$ws.Range("N94").Value = -3202

$ws.Range("H132").Value = 1043691
$ws.Range("I132").Value = 1073.7727
$ws.Range("J132").Value = 16335410
$ws.Range("K132").Value = 3221.3181
$ws.Range("L132").Value = 49006230
$ws.Range("M132").Value = -691.3181
$ws.Range("N132").Value = -49011290

$ws.Range("H137").Value = 2779547.5
$ws.Range("I137").Value = 3227353.5
$ws.Range("J137").Value = 3150
$ws.Range("K137").Value = 9682060.5
$ws.Range("L137").Value = 9450
$ws.Range("M137").Value = -9679510.5
$ws.Range("N137").Value = -14550

$ws.Range("H138").Value = 2350070
$ws.Range("I138").Value = 1339.8889
$ws.Range("J138").Value = 3791336.2
$ws.Range("K138").Value = 4019.6667
$ws.Range("L138").Value = 11374008.6
$ws.Range("M138").Value = 1120.3333
$ws.Range("N138").Value = -11384288.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14426.096
$ws.Range("I32").Value = 14961.197
$ws.Range("J32").Value = 11092
$ws.Range("K32").Value = 14961.197
$ws.Range("L32").Value = 11092
$ws.Range("M32").Value = -14674.197
$ws.Range("N32").Value = -11666

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H45").Value = 2044.4546
$ws.Range("I45").Value = 3300
$ws.Range("J45").Value = 1327
$ws.Range("K45").Value = 3300
$ws.Range("L45").Value = 1327
$ws.Range("M45").Value = -2923
$ws.Range("N45").Value = -2081

$ws.Range("H61").Value = 77078030
$ws.Range("I61").Value = 91000856
$ws.Range("J61").Value = 502500
$ws.Range("K61").Value = 91000856
$ws.Range("L61").Value = 502500
$ws.Range("M61").Value = -91000644
$ws.Range("N61").Value = -502924

$ws.Range("H74").Value = 6001946.5
$ws.Range("I74").Value = 7607590
$ws.Range("J74").Value = 114586.664
$ws.Range("K74").Value = 7607590
$ws.Range("L74").Value = 114586.664
$ws.Range("M74").Value = -7606716
$ws.Range("N74").Value = -116334.664

$ws.Range("H77").Value = 6001946.5
$ws.Range("I77").Value = 7607590
$ws.Range("J77").Value = 114586.664
$ws.Range("K77").Value = 38037950
$ws.Range("L77").Value = 572933.3200000001
$ws.Range("M77").Value = -38033582
$ws.Range("N77").Value = -581669.3200000001

$ws.Range("H132").Value = 47691.56
$ws.Range("I132").Value = 27098.395
$ws.Range("J132").Value = 204199.6
$ws.Range("K132").Value = 81295.185
$ws.Range("L132").Value = 612598.8
$ws.Range("M132").Value = -78765.185
$ws.Range("N132").Value = -617658.8

$ws.Range("H136").Value = 77078030
$ws.Range("I136").Value = 91000856
$ws.Range("J136").Value = 502500
$ws.Range("K136").Value = 273002568
$ws.Range("L136").Value = 1507500
$ws.Range("M136").Value = -273000018
$ws.Range("N136").Value = -1512600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1399.5454
$ws.Range("I20").Value = 1181.2
$ws.Range("J20").Value = 1581.5
$ws.Range("K20").Value = 1181.2
$ws.Range("L20").Value = 1581.5
$ws.Range("M20").Value = -934.2
$ws.Range("N20").Value = -2075.5

$ws.Range("H99").Value = 1217.55
$ws.Range("I99").Value = 1315.3334
$ws.Range("J99").Value = 924.2
$ws.Range("K99").Value = 1315.3334
$ws.Range("L99").Value = 924.2
$ws.Range("M99").Value = 182.6666
$ws.Range("N99").Value = -3920.2

$ws.Range("H134").Value = 2148.131
$ws.Range("I134").Value = 1651.7451
$ws.Range("J134").Value = 4679.7
$ws.Range("K134").Value = 4955.2353
$ws.Range("L134").Value = 14039.1
$ws.Range("M134").Value = -2420.2353
$ws.Range("N134").Value = -19109.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 25001070
$ws.Range("I58").Value = 26316600
$ws.Range("J58").Value = 6000.5
$ws.Range("K58").Value = 26316600
$ws.Range("L58").Value = 6000.5
$ws.Range("M58").Value = -26316397
$ws.Range("N58").Value = -6406.5

$ws.Range("H63").Value = 30271
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 30271
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 30271
$ws.Range("N63").Value = -31643

$ws.Range("H66").Value = 30271
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 30271
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 90813
$ws.Range("N66").Value = -97677

$ws.Range("H107").Value = 354.24243
$ws.Range("I107").Value = 336.81818
$ws.Range("J107").Value = 389.0909
$ws.Range("K107").Value = 336.81818
$ws.Range("L107").Value = 389.0909
$ws.Range("M107").Value = 1583.18182
$ws.Range("N107").Value = -4229.0909

$ws.Range("H132").Value = 18948.982
$ws.Range("I132").Value = 1093.4694
$ws.Range("J132").Value = 128314
$ws.Range("K132").Value = 3280.4082
$ws.Range("L132").Value = 384942
$ws.Range("M132").Value = -750.4081999999999
$ws.Range("N132").Value = -390002

$ws.Range("H136").Value = 25001070
$ws.Range("I136").Value = 26316600
$ws.Range("J136").Value = 6000.5
$ws.Range("K136").Value = 78949800
$ws.Range("L136").Value = 18001.5
$ws.Range("M136").Value = -78947250
$ws.Range("N136").Value = -23101.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 817.25
$ws.Range("I5").Value = 542.7143
$ws.Range("J5").Value = 930.2941
$ws.Range("K5").Value = 1628.1429
$ws.Range("L5").Value = 2790.8823
$ws.Range("M5").Value = -1516.1429
$ws.Range("N5").Value = -3014.8823

$ws.Range("H131").Value = 988.8046000000001
$ws.Range("I131").Value = 597.25
$ws.Range("J131").Value = 1028.4557
$ws.Range("K131").Value = 1791.75
$ws.Range("L131").Value = 3085.3671
$ws.Range("M131").Value = 3248.25
$ws.Range("N131").Value = -13165.3671

$ws.Range("H132").Value = 972.8125
$ws.Range("I132").Value = 687.7273
$ws.Range("J132").Value = 1600
$ws.Range("K132").Value = 6189.545700000001
$ws.Range("L132").Value = 14400
$ws.Range("M132").Value = -3659.545700000001
$ws.Range("N132").Value = -19460

$ws.Range("H135").Value = 817.25
$ws.Range("I135").Value = 542.7143
$ws.Range("J135").Value = 930.2941
$ws.Range("K135").Value = 4884.428699999999
$ws.Range("L135").Value = 8372.6469
$ws.Range("M135").Value = -2349.428699999999
$ws.Range("N135").Value = -13442.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 37865.582
$ws.Range("I70").Value = 58017.527
$ws.Range("J70").Value = 5958.3335
$ws.Range("K70").Value = 58017.527
$ws.Range("L70").Value = 5958.3335
$ws.Range("M70").Value = -57747.527
$ws.Range("N70").Value = -6498.3335

$ws.Range("H73").Value = 37865.582
$ws.Range("I73").Value = 58017.527
$ws.Range("J73").Value = 5958.3335
$ws.Range("K73").Value = 58017.527
$ws.Range("L73").Value = 5958.3335
$ws.Range("M73").Value = -57081.527
$ws.Range("N73").Value = -7830.3335

$ws.Range("H97").Value = 1241.75
$ws.Range("I97").Value = 1417.091
$ws.Range("J97").Value = 856
$ws.Range("K97").Value = 1417.091
$ws.Range("L97").Value = 856
$ws.Range("M97").Value = -921.0909999999999
$ws.Range("N97").Value = -1848

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3272.2942
$ws.Range("I16").Value = 1108.6
$ws.Range("J16").Value = 19500
$ws.Range("K16").Value = 1108.6
$ws.Range("L16").Value = 19500
$ws.Range("M16").Value = -938.5999999999999
$ws.Range("N16").Value = -19840

$ws.Range("H132").Value = 32032.42
$ws.Range("I132").Value = 12698.326
$ws.Range("J132").Value = 254374.5
$ws.Range("K132").Value = 38094.978
$ws.Range("L132").Value = 763123.5
$ws.Range("M132").Value = -35564.978
$ws.Range("N132").Value = -768183.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

$ws.Range("H45").Value = 13000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 13000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 13000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -13982

$ws.Range("H81").Value = 1980.2307
$ws.Range("I81").Value = 722.375
$ws.Range("J81").Value = 2539.2778
$ws.Range("K81").Value = 1444.75
$ws.Range("L81").Value = 5078.5556
$ws.Range("M81").Value = -383.75
$ws.Range("N81").Value = -7200.5556

$ws.Range("H84").Value = 1980.2307
$ws.Range("I84").Value = 722.375
$ws.Range("J84").Value = 2539.2778
$ws.Range("K84").Value = 7223.75
$ws.Range("L84").Value = 25392.778
$ws.Range("M84").Value = -1919.75
$ws.Range("N84").Value = -36000.778

$ws.Range("H107").Value = 301.75
$ws.Range("I107").Value = 268.85715
$ws.Range("J107").Value = 327.33334
$ws.Range("K107").Value = 806.5714499999999
$ws.Range("L107").Value = 982.0000200000001
$ws.Range("M107").Value = 1113.42855
$ws.Range("N107").Value = -4822.00002
